$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update input values (formulas recalc automatically)
$ws.Range("E3").Value = 70
$ws.Range("G3").Value = 185
$ws.Range("E4").Value = 73.498999999999995
$ws.Range("G4").Value = 262.33699999999999
$ws.Range("E7").Value = 3700

# Remove row 10 (L10 formula) entirely
$ws.Rows("10").Delete()

# Update selection to match new view
$ws.Range("F10").Select()
